$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update is_active (column D) to FALSE for the "Others" rows (OTH) in each language
$ws.Range("D4").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("D10").Value = $false

# Update the last active cell selection to D12
$ws.Range("D12").Select()
